$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 6 (2025) values per new data:
# total_customers 406 -> 407
$ws.Range("C6").Value = 407
# new_customers 99 -> 100
$ws.Range("E6").Value = 100
# new_rate 24.38423645320197 -> 24.57002457002457
$ws.Range("G6").Value = 24.57002457002457
# returning_rate 75.61576354679804 -> 75.42997542997543
$ws.Range("H6").Value = 75.42997542997543
